# Arbeitspakete_Liste_NEU.xlsx - add new work package row "Berechtigungskonzept
# erstellen" as the new AP #5 (old AP #5 "Kernsystem implementieren" and all
# subsequent rows shift down by one), and update the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 6 ("Kernsystem implementieren"),
# pushing rows 6..19 down to 7..20. Excel copies the row-6 formatting onto the
# freshly inserted row, which matches the target (s="4"/"6"/"8"/"3"/"3"/"1").
$ws.Rows.Item(6).Insert()

# Fill in the new work package row.
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Berechtigungskonzept erstellen"

# Renumber column A ("Nr.") for the shifted rows so the sequence stays
# consecutive (old row r held value r-1 for r in 7..19 -> now needs r-1 again
# after the shift, and the row that used to be blank at A7 gets filled too).
for ($r = 7; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Update the sheet's active selection to A2:A20 (anchor A2), matching the
# saved view state after the edit.
$ws.Range("A2:A20").Select()
